$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New item: "Branch of the Cursed Tree Barnstokkr" (replaces the old "piece of wood")
$ws.Range("A64").Value = 6610062
$ws.Range("B64").Value = "Branch of the Cursed Tree Barnstokkr"

$ws.Range("A65").Value = 6610063
$ws.Range("B65").Value = "This twisted length of wood, carved with strange symbols, is one of the branches of the legendary Barnstokkr, a cursed tree located south of Kokari Wilds, in lands unknown. The legend says that the tree is doomed to burn forever because of a curse invoked by a witch who, before dying, killed all the Templars that was hunting, causing a huge explosion of flames malignant, near the oldest tree the forest where she lived. Since then, the tree burns eternally in these flames that are full of hatred, bitterness and rage. Probably, the symbols that were carved on the branch, used to contain and control the terrible evil energy emanating from these flames. Do not want to know where your Mabari found this branch, much less whom he belonged."

# Match formatting of the existing, non-highlighted rows (plain banded fill, no wrap)
$ws.Range("A61:B61").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the new selection left behind by the author's edit
[void]$ws.Range("A64:A65").Select()
